# Refresh the crypto price ('Price') and volume change ('Volume(1h)') columns
# with the latest scraped figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ('Price') ------------------------------------------------
# These are plain-text cells in the source sheet (prices are pre-formatted
# strings, not numbers). Assigning a numeric-looking string to .Value would
# normally get auto-converted to a real number by Excel, so each cell is
# temporarily switched to Text format, written, then restored to the default
# 'Normal' style so no residual number formatting is left behind.
$priceUpdates = [ordered]@{
    D2 = '70.206.24'
    D3 = '3.605.60'
    D5 = '604.95'
    D6 = '196.08'
    D11 = '53.84'
    D13 = '9.58'
    D14 = '4.177.65'
    D15 = '13.11'
    D16 = '597.07'
    D17 = '70.334.48'
    D18 = '19.10'
    D19 = '3.609.24'
    D21 = '0.997'
    D22 = '17.78'
    D23 = '5.20'
    D24 = '102.11'
    D25 = '4.62'
    D27 = '10.75'
    D28 = '9.63'
    D29 = '33.79'
    D30 = '4.76'
    D31 = '7.15'
    D32 = '12.28'
    D34 = '63.28'
    D35 = '0.0₃0893'
    D36 = '3.910.05'
    D37 = '3.09'
    D39 = '519.13'
    D40 = '36.93'
    D42 = '3.53'
    D44 = '0.0453'
    D45 = '3.40'
    D46 = '2.85'
    D48 = '8.63'
}
foreach ($cell in $priceUpdates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cell]
    $range.Style = "Normal"
}

# --- Column E ('Volume(1h)') --------------------------------------------
# Already padded/percent strings (e.g. "  +0.34%  "), so they round-trip
# as text without any special handling.
$volumeUpdates = [ordered]@{
    E2 = '  +0.34%  '
    E3 = '  +2.21%  '
    E4 = '  +0.06%  '
    E5 = '  +0.20%  '
    E6 = '  -0.13%  '
    E7 = '  +0.04%  '
    E8 = '  +0.05%  '
    E9 = '  -2.04%  '
    E10 = '  -1.18%  '
    E11 = '  -0.63%  '
    E12 = '  +0.23%  '
    E13 = '  +0.21%  '
    E14 = '  +2.32%  '
    E15 = '  +3.62%  '
    E16 = '  -1.06%  '
    E17 = '  +0.37%  '
    E18 = '  -0.42%  '
    E19 = '  +2.47%  '
    E20 = '  +1.44%  '
    E21 = '  +0.00%  '
    E22 = '  -2.87%  '
    E23 = '  -1.15%  '
    E24 = '  -1.37%  '
    E25 = '  -0.02%  '
    E26 = '  -2.28%  '
    E27 = '  -2.01%  '
    E28 = '  -0.77%  '
    E29 = '  +0.23%  '
    E30 = '  +5.86%  '
    E32 = '  -3.78%  '
    E33 = '  +0.49%  '
    E34 = '  -0.23%  '
    E35 = '  +7.17%  '
    E36 = '  +4.29%  '
    E37 = '  -0.11%  '
    E38 = '  -0.01%  '
    E39 = '  +5.66%  '
    E40 = '  +0.08%  '
    E41 = '  -1.25%  '
    E42 = '  -2.15%  '
    E43 = '  -1.97%  '
    E44 = '  -1.03%  '
    E45 = '  +1.96%  '
    E46 = '  +0.67%  '
    E47 = '  -0.17%  '
    E48 = '  -0.62%  '
    E49 = '  -0.23%  '
    E50 = '  +2.23%  '
    E51 = '  +0.80%  '
}
foreach ($cell in $volumeUpdates.Keys) {
    $ws.Range($cell).Value = $volumeUpdates[$cell]
}

